{"js": "// Replace each three-digit-by-one-digit multiplication \"problem=answer\"\n// text with its new value. Every old value is unique in the document, so\n// a plain text search finds exactly one hit per pair.\nconst replacements = [\n  [\"374\u00d79=3366\", \"890\u00d79=8010\"],\n  [\"190\u00d77=1330\", \"837\u00d78=6696\"],\n  [\"681\u00d79=6129\", \"302\u00d73=906\"],\n  [\"260\u00d75=1300\", \"376\u00d78=3008\"],\n  [\"775\u00d76=4650\", \"666\u00d78=5328\"],\n  [\"224\u00d74=896\", \"140\u00d74=560\"],\n  [\"103\u00d77=721\", \"613\u00d78=4904\"],\n  [\"610\u00d72=1220\", \"729\u00d75=3645\"],\n  [\"459\u00d75=2295\", \"169\u00d79=1521\"],\n  [\"697\u00d78=5576\", \"883\u00d75=4415\"],\n  [\"437\u00d77=3059\", \"274\u00d74=1096\"],\n  [\"445\u00d79=4005\", \"830\u00d76=4980\"],\n  [\"264\u00d77=1848\", \"878\u00d73=2634\"],\n  [\"726\u00d79=6534\", \"605\u00d75=3025\"],\n  [\"550\u00d75=2750\", \"797\u00d77=5579\"],\n  [\"659\u00d72=1318\", \"991\u00d74=3964\"],\n  [\"844\u00d76=5064\", \"739\u00d72=1478\"],\n  [\"280\u00d76=1680\", \"285\u00d77=1995\"],\n  [\"766\u00d73=2298\", \"514\u00d75=2570\"],\n  [\"579\u00d75=2895\", \"821\u00d79=7389\"],\n  [\"689\u00d73=2067\", \"328\u00d77=2296\"],\n  [\"326\u00d73=978\", \"633\u00d78=5064\"],\n  [\"545\u00d78=4360\", \"425\u00d73=1275\"],\n  [\"200\u00d73=600\", \"819\u00d74=3276\"],\n  [\"950\u00d78=7600\", \"635\u00d76=3810\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each row holds the original \"problem=answer\" text and its replacement.\n# Every original value is unique in the document, so Find/Execute replaces\n# exactly one occurrence per pair.\n$pairs = @(\n  @(\"374\u00d79=3366\", \"890\u00d79=8010\"),\n  @(\"190\u00d77=1330\", \"837\u00d78=6696\"),\n  @(\"681\u00d79=6129\", \"302\u00d73=906\"),\n  @(\"260\u00d75=1300\", \"376\u00d78=3008\"),\n  @(\"775\u00d76=4650\", \"666\u00d78=5328\"),\n  @(\"224\u00d74=896\", \"140\u00d74=560\"),\n  @(\"103\u00d77=721\", \"613\u00d78=4904\"),\n  @(\"610\u00d72=1220\", \"729\u00d75=3645\"),\n  @(\"459\u00d75=2295\", \"169\u00d79=1521\"),\n  @(\"697\u00d78=5576\", \"883\u00d75=4415\"),\n  @(\"437\u00d77=3059\", \"274\u00d74=1096\"),\n  @(\"445\u00d79=4005\", \"830\u00d76=4980\"),\n  @(\"264\u00d77=1848\", \"878\u00d73=2634\"),\n  @(\"726\u00d79=6534\", \"605\u00d75=3025\"),\n  @(\"550\u00d75=2750\", \"797\u00d77=5579\"),\n  @(\"659\u00d72=1318\", \"991\u00d74=3964\"),\n  @(\"844\u00d76=5064\", \"739\u00d72=1478\"),\n  @(\"280\u00d76=1680\", \"285\u00d77=1995\"),\n  @(\"766\u00d73=2298\", \"514\u00d75=2570\"),\n  @(\"579\u00d75=2895\", \"821\u00d79=7389\"),\n  @(\"689\u00d73=2067\", \"328\u00d77=2296\"),\n  @(\"326\u00d73=978\", \"633\u00d78=5064\"),\n  @(\"545\u00d78=4360\", \"425\u00d73=1275\"),\n  @(\"200\u00d73=600\", \"819\u00d74=3276\"),\n  @(\"950\u00d78=7600\", \"635\u00d76=3810\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
